$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "BNS - Coalition Bulgarian People's Union (Koalitsiya Bulgarski Naroden Suyoz , BNS)"
$ws.Range("C1").Value = "DPS - Movement for Rights and Freedoms (Dvizehnie za Prava i Svobodi, DPS)"
$ws.Range("D1").Value = "DSB - Democrats for Strong Bulgaria (Demokrati za Silna Bulgaria, DSB)"
$ws.Range("E1").Value = "KA - Coalition Attack (Koalitsiya 'Ataka', KA)"
$ws.Range("F1").Value = "KB - Coalition for Bulgaria (Koalitsia za Bulgaria, KB)"
$ws.Range("G1").Value = "NDSV - National Movement for Stability and Progress (Nacionalno dviženie za stabilnost i văzhod, NDSV)"
$ws.Range("H1").Value = "ODS - United Democratic Forces (Obedineni Demokratichni Sili, ODS)"
$ws.Range("I1").Value = "GERB - Citizens for a European Development of Bulgaria (Grazhdani za Evropeisko Razvitie na Bulgaria, GERB)"
$ws.Range("J1").Value = "RZS - Order, Law and Justice (Red, Zakonnost I Spravedlivost, RZS)"
$ws.Range("K1").Value = "SK - Blue Coalition (Sinyata Koalitzia, SK)"
$ws.Range("L1").Value = "ABV - ABV-Alternative for Bulgarian Revival (Alternativa Za Bulgarsko Vuzrazhdane, ABV)"
$ws.Range("M1").Value = "BBT - Bulgaria without Censorship (Bulgaria Without Censorship-VMRO-BNU-Gergovden, BBT)"
$ws.Range("N1").Value = "PF - Patriotic Front (Patriotichen Front, PF)"
$ws.Range("O1").Value = "RB - Reformist Bloc (Reformatorski Blok, RB)"
$ws.Range("P1").Value = "OP - United Patriots  (Obedineni Patrioti, OP)"
$ws.Range("Q1").Value = "Will - Will (Volya, Will)"
